$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a cell to an explicit text value (preserves "General" number format / no style
# change) by using a leading quote-prefix, then resetting the style so no quotePrefix
# flag lingers on the cell -- matches the original inlineStr (text) cells exactly.
function Set-TextValue {
    param($cell, $value)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Updated symbol list on Fri Feb  3 03:55:41 UTC 2023 with GitHub Actions
Set-TextValue $ws.Range("D2") "321.37"
Set-TextValue $ws.Range("E2") "-1.78%"
Set-TextValue $ws.Range("D3") "39.38"
Set-TextValue $ws.Range("E3") "-0.87%"
Set-TextValue $ws.Range("D4") "5.914"
Set-TextValue $ws.Range("E4") "12.68%"
Set-TextValue $ws.Range("D5") "0.08013"
Set-TextValue $ws.Range("E5") "-0.98%"
Set-TextValue $ws.Range("D6") "4.560"
Set-TextValue $ws.Range("E6") "0.86%"
Set-TextValue $ws.Range("D7") "8.649"
Set-TextValue $ws.Range("E7") "0.17%"
Set-TextValue $ws.Range("D8") "1.930"
Set-TextValue $ws.Range("E8") "0.84%"
Set-TextValue $ws.Range("D9") "0.9340"
Set-TextValue $ws.Range("E9") "-0.02%"
Set-TextValue $ws.Range("D10") "0.1253"
Set-TextValue $ws.Range("E10") "-5.48%"
Set-TextValue $ws.Range("D11") "0.1959"
Set-TextValue $ws.Range("E11") "-0.12%"
Set-TextValue $ws.Range("D12") "8.774"
Set-TextValue $ws.Range("E12") "21.31%"
Set-TextValue $ws.Range("D13") "0.09107"
Set-TextValue $ws.Range("E13") "-0.54%"
Set-TextValue $ws.Range("D14") "0.03552"
Set-TextValue $ws.Range("E14") "3.15%"
Set-TextValue $ws.Range("D15") "0.09571"
Set-TextValue $ws.Range("E15") "0.16%"
Set-TextValue $ws.Range("D16") "0.001298"
Set-TextValue $ws.Range("E16") "-7.15%"
Set-TextValue $ws.Range("D17") "0.006249"
Set-TextValue $ws.Range("E17") "2.06%"
Set-TextValue $ws.Range("D18") "3.356"
Set-TextValue $ws.Range("E18") "-0.09%"
Set-TextValue $ws.Range("D19") "2.942"
Set-TextValue $ws.Range("E19") "-0.55%"
Set-TextValue $ws.Range("D20") "0.3536"
Set-TextValue $ws.Range("E20") "0.05%"
Set-TextValue $ws.Range("E21") "8.06%"
Set-TextValue $ws.Range("D22") "0.2410"
Set-TextValue $ws.Range("E22") "4.28%"
Set-TextValue $ws.Range("D23") "0.04445"
Set-TextValue $ws.Range("E23") "0.15%"
Set-TextValue $ws.Range("E24") "3.22%"
Set-TextValue $ws.Range("D25") "0.004393"
Set-TextValue $ws.Range("E25") "0.89%"
Set-TextValue $ws.Range("E26") "-11.64%"
Set-TextValue $ws.Range("E27") "0.05%"
Set-TextValue $ws.Range("D39") "0.02404"
Set-TextValue $ws.Range("E39") "-3.30%"
Set-TextValue $ws.Range("D40") "0.05177"
Set-TextValue $ws.Range("E40") "-1.11%"
Set-TextValue $ws.Range("D41") "0.007437"
Set-TextValue $ws.Range("E41") "-3.25%"
Set-TextValue $ws.Range("D42") "0.009344"
Set-TextValue $ws.Range("E42") "8.21%"
Set-TextValue $ws.Range("D43") "0.1405"
Set-TextValue $ws.Range("E43") "-1.98%"
Set-TextValue $ws.Range("E44") "0.46%"
Set-TextValue $ws.Range("D45") "0.01124"
Set-TextValue $ws.Range("E45") "37.55%"
Set-TextValue $ws.Range("D46") "0.00006734"
Set-TextValue $ws.Range("E46") "0.96%"
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.04%"
Set-TextValue $ws.Range("E48") "5.44%"
Set-TextValue $ws.Range("E50") "0.04%"
Set-TextValue $ws.Range("E51") "0.04%"
